# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (used by the notes master)
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet" (used by the slide
#                            master, i.e. the theme that actually drives
#                            every slide's look)
# The authored change swaps the two themes' colour schemes, so the slide
# master (and therefore every slide) goes from the "Integral" / Red Violet
# palette back to the default Office colour palette.
#
# The PowerPoint object model doesn't give us byte-level access to the OOXML
# theme parts, but ColorScheme.Colors(i).RGB on the slide master's colour
# scheme writes straight through to the underlying <a:clrScheme> that the
# slides use, so we drive the swap through that.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.ColorScheme

# Target palette = the classic default "Office Theme" colour scheme
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), expressed as the packed
# BGR-in-a-long values PowerPoint's RGB property expects (R | G<<8 | B<<16).
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $scheme.Colors($i).RGB = $officeColors[$i - 1]
}
